$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13236
$ws1.Range("F5").Value = 129
$ws1.Range("F11").Value = 13182
$ws1.Range("F12").Value = 323
$ws1.Range("F26").Value = 196

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 13236
$ws4.Range("F6").Value = 129
$ws4.Range("F12").Value = 13182
$ws4.Range("F13").Value = 323
$ws4.Range("F29").Value = 196
